$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.307.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.695.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.98%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "498.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.99%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -5.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.705.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.80%  "
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("E11").Value = "  -3.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.345"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.60%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.170.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.399.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.703.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.75%  "
$ws.Range("E18").Value = "  -5.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.61%  "
$ws.Range("E20").Value = "  -6.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "332.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.33%  "
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.422"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.170"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0816"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.43%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.22%  "
$ws.Range("E37").Value = "  -6.75%  "
$ws.Range("E38").Value = "  -8.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.61%  "
$ws.Range("E40").Value = "  -4.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.175.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.24%  "
$ws.Range("E42").Value = "  -8.61%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0553"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.592"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.98%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.37%  "
$ws.Range("E48").Value = "  -4.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.79%  "
$ws.Range("E50").Value = "  -5.54%  "
$ws.Range("E51").Value = "  -3.70%  "
